# Edit corresponds to commit "#5: fund, bonds, otherbonds, antique done"
#
# Sheet4 "基金受益憑證" (fund): add a proper header row (row 1) with column
# labels, and append the standard metadata columns
# (property_category, category, date, legislator_name, legislator_id,
# source_file, index) to every data row, matching the pattern already used
# on sheet1/2/3.
#
# Sheet5 "具有相當價值之財產" (otherbonds/antique): same treatment - proper
# header row and the same trailing metadata columns appended.
#
# Sheet6 "保險" is unaffected in content (its shared-string indices merely
# shift because of new strings inserted earlier in the shared string table;
# Excel/iron_native manage that automatically when we simply leave the sheet
# untouched).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet4: 基金受益憑證 (fund)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$header4 = @("name","owner","dealer","quantity","face_value","currency","total","property_category","category","date","legislator_name","legislator_id","source_file","index")
for ($i = 0; $i -lt $header4.Length; $i++) {
    $ws4.Cells.Item(1, $i + 2).Value = $header4[$i]
}

$rows4 = @(
    @(2, "匯豐拉美", "蔡麗卿", "合作金庫", 20000, 10, "新臺幣", 200000, "fund", "normal", "'2011-12-22", "徐耀昌", 921, "tmpd3a41", 72),
    @(3, "聯博全高T", "蔡麗卿", "合作金庫", 1157.4, 4, "美金", 152500, "fund", "normal", "'2011-12-22", "徐耀昌", 921, "tmpd3a41", 73),
    @(4, "富蘭克林坦公司債", "蔡麗卿", "合作金庫", 764.526, 6.5399, "美金", 152500, "fund", "normal", "'2011-12-22", "徐耀昌", 921, "tmpd3a41", 74),
    @(5, "富蘭克林亞洲成長", "蔡麗卿", "合作金庫", 107.428, 29.7873, "美金", 97600, "fund", "normal", "'2011-12-22", "徐耀昌", 921, "tmpd3a41", 75),
    @(6, "富蘭克林天資美", "蔡麗卿", "合作金庫", 327.807, 9.7618, "美金", 97600, "fund", "normal", "'2011-12-22", "徐耀昌", 921, "tmpd3a41", 77)
)

foreach ($row in $rows4) {
    $r = $row[0]
    for ($c = 1; $c -lt $row.Length; $c++) {
        $ws4.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# ---------------------------------------------------------------------
# Sheet5: 具有相當價值之財產 (otherbonds)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$header5 = @("name","quantity","owner","total","property_category","category","date","legislator_name","legislator_id","source_file","index")
for ($i = 0; $i -lt $header5.Length; $i++) {
    $ws5.Cells.Item(1, $i + 2).Value = $header5[$i]
}

# Row 2 data (A2/C2/D2/E2 already correct; fill in the new trailing columns)
$ws5.Cells.Item(2, 6).Value = "otherbonds"
$ws5.Cells.Item(2, 7).Value = "normal"
$ws5.Cells.Item(2, 8).Value = "'2011-12-22"
$ws5.Cells.Item(2, 9).Value = "徐耀昌"
$ws5.Cells.Item(2, 10).Value = 921
$ws5.Cells.Item(2, 11).Value = "tmpd3a41"
$ws5.Cells.Item(2, 12).Value = 86
